# update api & mobile testing
# Row 5 (TC001 / positive) test data changed:
#   fullname (C5): "Archie" -> "Sa"
#   email    (D5): "archen@gmail.com" -> "sal@gmail.com"
# (the D5 hyperlink relationship itself is left untouched - only the
# displayed text/shared-string changes, matching the source diff which
# shows no change to the worksheet's .rels file)

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("C5").Value = "Sa"
$ws.Range("D5").Value = "sal@gmail.com"

# Leave the cursor where the author left it after the edit.
$ws.Range("G10").Select()
